# DesignFirst project save: update rule R30 "From" threshold value (cell C10)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 100
